# The workbook has a small results table with headers in row 1:
#   B1:Feature  C1:DoF  D1:P Value  E1:Chi Square  F1:Is Significant
# and two data rows (2 and 3) below it.
#
# This change inserts two new columns ("Observed" and "Expected") between
# "Chi Square" and "Is Significant", shifting the old F column (Is
# Significant) to H, and fills the new F/G data cells with the
# observed/expected contingency-table strings for both data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at F:G - this shifts the existing F column
# (and its formatting) to H, preserving the header cell style.
$ws.Columns("F:G").Insert()

# New header cells for the inserted columns.
$ws.Range("F1").Value = "Observed"
$ws.Range("G1").Value = "Expected"

# New data values for the two result rows.
$ws.Range("F2").Value = "[ 86 159] ; [15  5]"
$ws.Range("G2").Value = "[ 93.37735849 151.62264151] ; [ 7.62264151 12.37735849]"

$ws.Range("F3").Value = "[ 86 159] ; [15  5]"
$ws.Range("G3").Value = "[ 93.37735849 151.62264151] ; [ 7.62264151 12.37735849]"
